$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -8
$ws.Range("F4").Value = 7
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 1
